# "add hw10 & group gits"
# Add each group's GitHub repo URL (git_page, col E) and final-project
# filename (file, col F) next to the existing roster data, then tidy up
# the sheet (drop the stray fill formatting on cols A/B/D, autosize the
# new-ish hyperlink column C, set the page to portrait, and leave the
# selection where the last edit happened).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- group -> (git_page url, file name) --------------------------------
$groupInfo = @{
    "a" = @("https://github.com/seul-b/edld-final",            "in-progress.Rmd")
    "b" = @("https://github.com/emaduneme/EDLD_651_Ghana",      "Main Markdown.Rmd")
    "c" = @("https://github.com/tianwalker44/EDLD_Final",       "Final_Groupof5.Rmd")
    "d" = @("https://github.com/haithamanbar/Oregon-made",      "Final Project.Rmd")
}

# --- clear the old "applyFill" styling Excel had stamped on cols A/B/D -
# (the hyperlink formatting in col C is left untouched)
$ws.Range("A2:B17").ClearFormats()
$ws.Range("D2:D17").ClearFormats()

# --- fill in git_page / file for every roster row (2-19) ---------------
for ($r = 2; $r -le 19; $r++) {
    $grp = $ws.Cells.Item($r, 4).Value2
    $info = $groupInfo[$grp]
    if ($info) {
        $ws.Cells.Item($r, 5).Value = $info[0]
        $ws.Cells.Item($r, 6).Value = $info[1]
    }
}

# --- autosize the hyperlink column now that data has changed -----------
$ws.Columns("C").EntireColumn.AutoFit()

# --- print setup: portrait -------------------------------------------
$ws.PageSetup.Orientation = 1

# --- leave selection on the last cell touched ---------------------------
$ws.Range("E15").Select()
